$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 58564
$ws.Range("B3").Value = 57805
$ws.Range("B4").Value = 58041
$ws.Range("B5").Value = 58321
$ws.Range("B6").Value = 58256
$ws.Range("B7").Value = 57881
